$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from B1 to C1:D1, then set header text
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)
$ws.Range('C1').Value = 'Authors'
$ws.Range('D1').Value = 'Dates'

# Fill Authors (C) and Dates (D) columns for data rows
$ws.Range('C2').Value = '美妙歌声'
$ws.Range('D2').Value = '09-13 04:21'
$ws.Range('C3').Value = 'chaoshou22'
$ws.Range('D3').Value = '08-31 05:58'
$ws.Range('C4').Value = '摆渡人88'
$ws.Range('D4').Value = '08-30 04:01'
$ws.Range('C5').Value = '流风回血'
$ws.Range('D5').Value = '08-30 03:29'
$ws.Range('C6').Value = 'chaoshou22'
$ws.Range('D6').Value = '08-30 09:51'
$ws.Range('C7').Value = '奥力给888'
$ws.Range('D7').Value = '08-30 07:38'
$ws.Range('C8').Value = '摆渡人88'
$ws.Range('D8').Value = '08-29 03:40'
$ws.Range('C9').Value = '流风回血'
$ws.Range('D9').Value = '08-29 08:04'
$ws.Range('C10').Value = '股海小歧'
$ws.Range('D10').Value = '08-29 07:43'
$ws.Range('C11').Value = '股海小歧'
$ws.Range('D11').Value = '08-28 10:50'
$ws.Range('C12').Value = '我想找个男朋友'
$ws.Range('D12').Value = '08-26 02:12'
$ws.Range('C13').Value = '全景网'
$ws.Range('D13').Value = '08-22 04:42'
$ws.Range('C14').Value = 'BT财经'
$ws.Range('D14').Value = '08-22 03:00'
$ws.Range('C15').Value = '股友0R087A9259'
$ws.Range('D15').Value = '08-16 10:00'
$ws.Range('C16').Value = '股友0R087A9259'
$ws.Range('D16').Value = '08-16 09:39'
$ws.Range('C17').Value = '股友0R087A9259'
$ws.Range('D17').Value = '08-15 01:13'
$ws.Range('C18').Value = '智通财经'
$ws.Range('D18').Value = '08-14 06:30'
$ws.Range('C19').Value = '托塔金天王'
$ws.Range('D19').Value = '08-03 05:50'
$ws.Range('C20').Value = '十优缠论'
$ws.Range('D20').Value = '07-31 08:27'
$ws.Range('C21').Value = '祥云0303'
$ws.Range('D21').Value = '07-30 12:16'
$ws.Range('C22').Value = 'GUyou868C06F186'
$ws.Range('D22').Value = '07-28 11:11'
$ws.Range('C23').Value = 'GUyou10K3531K59'
$ws.Range('D23').Value = '07-28 10:36'
$ws.Range('C24').Value = 'GU181A7006B9948'
$ws.Range('D24').Value = '07-28 10:26'
$ws.Range('C25').Value = '白云黄鹤'
$ws.Range('D25').Value = '07-21 10:10'
$ws.Range('C26').Value = '格隆汇'
$ws.Range('D26').Value = '07-13 06:15'
$ws.Range('C27').Value = '格隆汇'
$ws.Range('D27').Value = '07-13 06:15'
$ws.Range('C28').Value = '格隆汇'
$ws.Range('D28').Value = '07-13 06:14'
$ws.Range('C29').Value = '格隆汇'
$ws.Range('D29').Value = '07-13 06:14'
$ws.Range('C30').Value = '格隆汇'
$ws.Range('D30').Value = '07-13 06:10'
$ws.Range('C31').Value = '祥云0303'
$ws.Range('D31').Value = '06-27 09:04'
$ws.Range('C32').Value = '函股论道'
$ws.Range('D32').Value = '06-15 08:38'
$ws.Range('C33').Value = '小散的成长路'
$ws.Range('D33').Value = '06-06 07:29'
$ws.Range('C34').Value = '论股浅水'
$ws.Range('D34').Value = '06-05 06:20'
$ws.Range('C35').Value = '祥云0303'
$ws.Range('D35').Value = '06-02 10:05'
$ws.Range('C36').Value = '股友10g83135g7'
$ws.Range('D36').Value = '05-31 01:26'
$ws.Range('C37').Value = '函股论道'
$ws.Range('D37').Value = '05-30 11:31'
$ws.Range('C38').Value = '有魅力之霹雳火'
$ws.Range('D38').Value = '05-30 09:55'
$ws.Range('C39').Value = '函股论道'
$ws.Range('D39').Value = '05-30 12:49'
$ws.Range('C40').Value = '账户已注销'
$ws.Range('D40').Value = '05-29 09:57'
$ws.Range('C41').Value = '李百万1'
$ws.Range('D41').Value = '05-29 03:02'
$ws.Range('C42').Value = '股势仁升'
$ws.Range('D42').Value = '05-29 04:36'
$ws.Range('C43').Value = '小懒财富日记'
$ws.Range('D43').Value = '05-26 06:39'
$ws.Range('C44').Value = '函股论道'
$ws.Range('D44').Value = '05-26 07:23'
$ws.Range('C45').Value = '股友8229A1'
$ws.Range('D45').Value = '05-26 10:47'
$ws.Range('C46').Value = '李百万1'
$ws.Range('D46').Value = '05-25 10:12'
$ws.Range('C47').Value = '腾飞老马'
$ws.Range('D47').Value = '05-25 08:24'
$ws.Range('C48').Value = '宝宝铁粉'
$ws.Range('D48').Value = '05-25 08:27'
$ws.Range('C49').Value = '函股论道'
$ws.Range('D49').Value = '05-25 05:06'
$ws.Range('C50').Value = '流浪股民77'
$ws.Range('D50').Value = '05-24 05:14'
$ws.Range('C51').Value = '看究竟app'
$ws.Range('D51').Value = '05-24 04:18'
$ws.Range('C52').Value = '财经小花姐'
$ws.Range('D52').Value = '05-23 07:21'
$ws.Range('C53').Value = '山科智能'
$ws.Range('D53').Value = '04-20 11:28'
$ws.Range('C54').Value = '智通财经'
$ws.Range('D54').Value = '03-27 06:18'
$ws.Range('C55').Value = '祥云0303'
$ws.Range('D55').Value = '03-28 12:07'
$ws.Range('C56').Value = '广888888发'
$ws.Range('D56').Value = '03-01 11:57'
$ws.Range('C57').Value = '智研咨询'
$ws.Range('D57').Value = '02-27 10:02'
$ws.Range('C58').Value = '番茄小号'
$ws.Range('D58').Value = '01-13 02:33'
$ws.Range('C59').Value = '泡财经APP'
$ws.Range('D59').Value = '11-02 06:46'
$ws.Range('C60').Value = '哈哈哈哈哈哈711'
$ws.Range('D60').Value = '08-27 02:29'
$ws.Range('C61').Value = '大吉为富de火舞你'
$ws.Range('D61').Value = '08-16 06:50'
$ws.Range('C62').Value = '大海微蓝'
$ws.Range('D62').Value = '07-21 11:40'
$ws.Range('C63').Value = '老金天下'
$ws.Range('D63').Value = '07-21 08:58'
$ws.Range('C64').Value = 'A财神到A'
$ws.Range('D64').Value = '06-10 09:30'
$ws.Range('C65').Value = '大司马论市'
$ws.Range('D65').Value = '05-17 10:50'
$ws.Range('C66').Value = '指数增强小霸王'
$ws.Range('D66').Value = '04-29 12:24'
$ws.Range('C67').Value = '山科智能'
$ws.Range('D67').Value = '04-23 12:15'
$ws.Range('C68').Value = '稳操胜券bqn8jr'
$ws.Range('D68').Value = '04-18 08:57'
$ws.Range('C69').Value = '高尔础'
$ws.Range('D69').Value = '02-16 09:41'
$ws.Range('C70').Value = '蔚蓝大海A'
$ws.Range('D70').Value = '10-25 04:38'
$ws.Range('C71').Value = 'guy1234'
$ws.Range('D71').Value = '10-08 10:28'
$ws.Range('C72').Value = 'guy1234'
$ws.Range('D72').Value = '10-07 07:21'
$ws.Range('C73').Value = 'kk11000'
$ws.Range('D73').Value = '10-05 11:48'
$ws.Range('C74').Value = '蔚蓝大海A'
$ws.Range('D74').Value = '08-10 07:24'
$ws.Range('C75').Value = '你是股神123456'
$ws.Range('D75').Value = '07-31 07:35'
$ws.Range('C76').Value = '喝茶喝茶'
$ws.Range('D76').Value = '04-24 05:21'
